$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "NOMBRES"
$ws.Range("C5").Value = "CI"
$ws.Range("D5").Value = "FECHA NACIMIENTO"
$ws.Range("E5").Value = "No LIBRO"
$ws.Range("F5").Value = "No PARTIDA"
$ws.Range("G5").Value = "USUARIO"

$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("G6").Select()
